$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters (B, C, E, F, G, H, I, J, L, M) mapped to their 1-based column indices.
$colIndexes = @(2, 3, 5, 6, 7, 8, 9, 10, 12, 13)

# New values for rows 2..25 (data rows 0..23), one row per array, values ordered
# to match $colIndexes (B, C, E, F, G, H, I, J, L, M).
$data = @(
    @(16.92146497504947, 9.771865390017005, 14.46341312524095, 39.81084986156147, 40.19851507583203, 17.06117069584081, 26.68969260311969, 8.358230682090989, 12.42906579138839, 17.09441340903945),
    @(16.45607380828435, 9.426491263954286, 14.48117275068085, 39.87480729873754, 40.23378982543115, 17.12149757096481, 26.81967784372227, 8.369292699519313, 12.42644612849389, 16.99401315226216),
    @(16.16612181287106, 9.206000172148705, 14.49283049462138, 39.92440456428495, 40.27103829153459, 17.16234439012146, 26.90585607835353, 8.376500638187355, 12.42618740715412, 16.93442781384993),
    @(16.04708585061116, 9.114106827968559, 14.49777097651256, 39.94720563271296, 40.29011865735435, 17.17994416678281, 26.94257116650946, 8.379542872252413, 12.42642269068722, 16.91068248447848),
    @(16.02727238185542, 9.098727178539145, 14.49860281970031, 39.9511479234304, 40.29352195488482, 17.1829241590801, 26.94876399664423, 8.380054382262674, 12.42648237411183, 16.90677248523897),
    @(16.16451976525532, 9.204769021851749, 14.49289635436112, 39.92470159167463, 40.2712798471024, 17.16257788625941, 26.90634477126809, 8.376541241425276, 12.42618919906535, 16.9341053818833),
    @(16.76197329630478, 9.654574272963126, 14.46938057033297, 39.8307550279667, 40.20743182933159, 17.08117991554309, 26.73318754748136, 8.361958791471061, 12.42788319478414, 17.05937866857691),
    @(17.89298438362104, 10.46675680615104, 14.42922323385452, 39.72874242211491, 40.20661143835579, 16.95187325894655, 26.44433373877166, 8.336644955313286, 12.44185537689861, 17.32057136093391),
    @(18.69026948734716, 11.01742699398651, 14.40332420837453, 39.70423011656637, 40.28255016698438, 16.87550875202876, 26.26328605602332, 8.320024561754629, 12.45853042349428, 17.52080188116627),
    @(19.04401517821062, 11.25736434825678, 14.3923190937536, 39.70407521754321, 40.33379365193656, 16.84484772808415, 26.18774877202049, 8.31288804398349, 12.46748866683185, 17.61344873368512),
    @(19.1765619930124, 11.34666308129232, 14.38826296421196, 39.70559935963966, 40.35559964608673, 16.83382610634395, 26.16013069138521, 8.31024624798075, 12.47107646909113, 17.64873424697341),
    @(19.14808016593451, 11.32750099126505, 14.38913158173051, 39.70520070047947, 40.35079658271358, 16.83617356998105, 26.1660347918598, 8.310812514169351, 12.47029510727861, 17.64112624714231),
    @(19.05494878622434, 11.26474249519721, 14.39198316575735, 39.70416888316148, 40.33553953452063, 16.84392915222069, 26.18545681993276, 8.312669488335585, 12.46777993028094, 17.61634777120288),
    @(18.99771618134863, 11.22609678128909, 14.39374432230241, 39.70374301907112, 40.32650679201018, 16.84875645947508, 26.1974819674489, 8.313814826563716, 12.46626471548086, 17.60119589993723),
    @(18.66696063599345, 11.00153028113372, 14.40405900997639, 39.70446166994677, 40.27953738763196, 16.87759480532443, 26.26836032110252, 8.320499455726967, 12.45797244456793, 17.51477670205148),
    @(18.46167104029, 10.86102891544264, 14.41058533446471, 39.70772037722293, 40.25500128204214, 16.89633250295796, 26.31359293101393, 8.324708653311855, 12.45323567444467, 17.46214535035846),
    @(18.3427587094744, 10.77922385813842, 14.41441221285535, 39.71062966393249, 40.24246075816544, 16.90749357549893, 26.34025135567385, 8.327169619053043, 12.45064053138, 17.43202147169692),
    @(18.30235753091949, 10.75135697410588, 14.41572049659363, 39.71179238562909, 40.23848468206775, 16.91133833083848, 26.34938751569638, 8.328009731336881, 12.44978412738046, 17.42184818668093),
    @(18.48361183421866, 10.87608854980938, 14.40988303240106, 39.70726635173902, 40.25745048470546, 16.89429811748336, 26.30871137835364, 8.32425644524362, 12.45372654070414, 17.46773286157456),
    @(19.08234289261689, 11.28321883288567, 14.3911425696853, 39.70442898951057, 40.3399557509699, 16.84163514223206, 26.17972529112534, 8.312122406786752, 12.4685134081077, 17.62362050271747),
    @(19.46538652464196, 11.54019202786032, 14.37954297489324, 39.71180024502081, 40.40787099474373, 16.81065142278234, 26.10117649583067, 8.304545457508732, 12.47931612276073, 17.72666992150659),
    @(19.26174254047719, 11.40388628223585, 14.3856746972595, 39.70702170945421, 40.37034397323099, 16.82687288541585, 26.14257150660696, 8.308557200889057, 12.4734469705811, 17.67157108387501),
    @(18.4736951596199, 10.8692832910158, 14.41020031014179, 39.70746839022598, 40.25633832332399, 16.89521665307317, 26.3109162926942, 8.324460760619761, 12.45350422094981, 17.46520632681904),
    @(17.59232895002962, 10.25490576988951, 14.43945192261792, 39.74750079330915, 40.19342934519574, 16.98359233530668, 26.51702106738538, 8.34314404886628, 12.4369435396042, 17.24836412977502)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $excelRow = $r + 2
    for ($c = 0; $c -lt $colIndexes.Length; $c++) {
        $ws.Cells.Item($excelRow, $colIndexes[$c]).Value = $rowValues[$c]
    }
}
